# Update to use context manager, file ready for merge
# Populate the budget rows (2-7) below the existing header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Income",  "w",    "`$123.00"),
    @("Income",  "2we2", "`$2.00"),
    @("Expense", "22",   "`$222.00"),
    @("Expense", "22",   "`$22.00"),
    @("Income",  "22",   "`$22.00"),
    @("Income",  "22",   "`$100.00")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    # Leading apostrophe forces these to stay plain text instead of being
    # auto-coerced into numbers / currency values by Excel's input parser.
    $ws.Cells.Item($r, 2).Value = "'" + $row[1]
    $ws.Cells.Item($r, 3).Value = "'" + $row[2]
    $r = $r + 1
}
